# This script applies a row-wise rotation of the batch-level figures
# (Batch No., Sale Rate, Quantity, Value columns) for several groups of
# duplicate stock rows in the report. Within each group of rows (which
# all share the same item in column C), the values in columns B, E, F
# and G are cyclically shifted by one row, while columns A, C and D are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(B, E, F, G) new values
$updates = @{
    136 = @(63902, 34.04, 2, 64.04000000000001)
    137 = @(48654, 38.26, -1, -32.02)

    146 = @(64350, 70.63, 2, 132.88)
    147 = @(57756, 79.37, -100, -6644)
    148 = @(53925, 79.37, 1, 66.44)

    246 = @(48706, 39.8, -144, -4795.2)
    247 = @(64973, 35.4, 81, 2697.3)

    292 = @(55373, 163.62, -94, -13562.32)
    293 = @(63520, 153.4, 75, 10821)
    294 = @(57802, 162.71, -79, -11334.92)
    296 = @(63571, 152.53, 6, 860.88)

    299 = @(55356, 54.04, -158, -7527.12)
    300 = @(63510, 50.66, 147, 7003.08)

    315 = @(63560, 134.87, 1, 126.86)
    316 = @(60325, 151.57, -102, -12939.72)

    420 = @(47097, 134.16, 15, 1684.2)
    421 = @(58047, 126.1, 42, 4432.68)

    465 = @(65069, 14.3, 2, 26.9)
    466 = @(53757, 16.08, -159, -2138.55)

    472 = @(64915, 20.98, 0, 0)
    473 = @(45695, 23.58, -36, -710.28)

    476 = @(64922, 20.98, 126, 2485.98)
    477 = @(45706, 23.58, -202, -3985.46)

    479 = @(45718, 19.38, -294, -4768.68)
    480 = @(64927, 17.26, 213, 3454.86)

    485 = @(64925, 13.97, 214, 2814.1)
    486 = @(45709, 15.69, -300, -3945)

    564 = @(64810, 291.22, 6, 1643.52)
    565 = @(53319, 310.64, -6, -1643.52)

    596 = @(64830, 34.9, 113, 3709.79)
    597 = @(60022, 37.22, -113, -3709.79)

    732 = @(65362, 43.44, 62, 2533.94)
    733 = @(65079, 43.44, 21, 858.27)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # Column B
    $ws.Cells.Item($row, 5).Value = $vals[1]   # Column E
    $ws.Cells.Item($row, 6).Value = $vals[2]   # Column F
    $ws.Cells.Item($row, 7).Value = $vals[3]   # Column G
}
